# Weekly update: insert a new price record as row 62, pushing the
# existing rows 62-117 down to 63-118 (dimension grows from R117 to R118).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 62 - shifts rows 62:117 down to 63:118 and
# extends the used range to A1:R118, carrying formatting (e.g. the date
# style on column D) down from the row above, same as Excel's UI insert.
$ws.Rows.Item(62).Insert()

# Fill in the new row 62 with the new weekly record. Columns A, B, C, E,
# F, G, H, I and R are constant for every row in this block, so reuse the
# same values as the rest of the table.
$ws.Range("A62").Value = 11
$ws.Range("B62").Value = "Vega Monumental Concepción"
$ws.Range("C62").Value = "Bíobío"
$ws.Range("D62").Value = 44658
$ws.Range("E62").Value = 8
$ws.Range("F62").Value = 100112043
$ws.Range("G62").Value = "Pepino ensalada"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 220
$ws.Range("K62").Value = 14000
$ws.Range("L62").Value = 15000
$ws.Range("M62").Value = 14545
$ws.Range("N62").Value = "$/caja 70 unidades"
$ws.Range("O62").Value = "Región del Maule"
$ws.Range("P62").Value = 208
$ws.Range("Q62").Value = 70
$ws.Range("R62").Value = "Hortaliza"
